# Add a new "Login" worksheet at the end of the workbook (after "Products"),
# populate it with a small URL/Username/Password table, box-border it, and
# hyperlink the URL + Password cells -- matching the authoring commit
# "My 13th commit through Intellij".

$wb = $excel.ActiveWorkbook

# --- 1. Create the new sheet after the last existing sheet ------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$loginSheet = $wb.Worksheets.Add($null, $lastSheet)
$loginSheet.Name = "Login"

# --- 2. Header row ------------------------------------------------------------
$loginSheet.Range("A1").Value = "URL"
$loginSheet.Range("B1").Value = "Username"
$loginSheet.Range("C1").Value = "Password"

# --- 3. Data rows (same credentials repeated on rows 2-4) ---------------------
$url = "http://49.249.28.218:8098/"
$user = "rmgyantra"
$pass = "rmgy@9999"

for ($r = 2; $r -le 4; $r++) {
    $loginSheet.Range("A$r").Value = $url
    $loginSheet.Range("B$r").Value = $user
    $loginSheet.Range("C$r").Value = $pass
}

# --- 4. Thin box border around the whole table ---------------------------------
$tableRange = $loginSheet.Range("A1:C4")
$tableRange.Borders.LineStyle = 1
$tableRange.Borders.Weight = 2

# --- 5. Column widths (closest achievable match to authored widths) -----------
$loginSheet.Columns.Item(1).ColumnWidth = 23
$loginSheet.Columns.Item(2).ColumnWidth = 8.5
$loginSheet.Columns.Item(3).ColumnWidth = 9.8

# --- 6. Hyperlinks on the URL and Password columns for every data row ---------
for ($r = 2; $r -le 4; $r++) {
    [void]$loginSheet.Hyperlinks.Add($loginSheet.Range("A$r"), $url)
    [void]$loginSheet.Hyperlinks.Add($loginSheet.Range("C$r"), "mailto:$pass")
}

# --- 7. Make "Login" the active sheet / selection, matching the saved state ---
[void]$loginSheet.Range("E5").Select()
